$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2690
$ws.Range("J40").Value = 2722.2222
$ws.Range("L40").Value = 2722.2222
$ws.Range("N40").Value = -3072.2222
$ws.Range("H48").Value = 550
$ws.Range("I48").Value = 550
$ws.Range("K48").Value = 1650
$ws.Range("M48").Value = -1358
$ws.Range("H56").Value = 550
$ws.Range("I56").Value = 550
$ws.Range("K56").Value = 1650
$ws.Range("M56").Value = -1116
$ws.Range("H70").Value = 1021657.6
$ws.Range("I70").Value = 2551769.8
$ws.Range("J70").Value = 1582.8334
$ws.Range("K70").Value = 7655309.399999999
$ws.Range("L70").Value = 4748.5002
$ws.Range("M70").Value = -7655039.399999999
$ws.Range("N70").Value = -5288.5002
$ws.Range("H73").Value = 1021657.6
$ws.Range("I73").Value = 2551769.8
$ws.Range("J73").Value = 1582.8334
$ws.Range("K73").Value = 7655309.399999999
$ws.Range("L73").Value = 4748.5002
$ws.Range("M73").Value = -7654373.399999999
$ws.Range("N73").Value = -6620.5002
$ws.Range("H112").Value = 1761.7894
$ws.Range("J112").Value = 1865.4
$ws.Range("L112").Value = 5596.200000000001
$ws.Range("N112").Value = -7812.200000000001
$ws.Range("H113").Value = 90913384
$ws.Range("I113").Value = 250002800
$ws.Range("J113").Value = 5143
$ws.Range("K113").Value = 250002800
$ws.Range("L113").Value = 5143
$ws.Range("M113").Value = -249999546
$ws.Range("N113").Value = -11651
$ws.Range("H137").Value = 1537.0358
$ws.Range("I137").Value = 1507.28
$ws.Range("J137").Value = 1785
$ws.Range("K137").Value = 4521.84
$ws.Range("L137").Value = 5355
$ws.Range("M137").Value = -1971.84
$ws.Range("N137").Value = -10455
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = $null

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 25000
$ws.Range("J27").Value = 25000
$ws.Range("L27").Value = 25000
$ws.Range("N27").Value = -25368
$ws.Range("H32").Value = 3647.9546
$ws.Range("I32").Value = 3702.6191
$ws.Range("K32").Value = 3702.6191
$ws.Range("M32").Value = -3415.6191
$ws.Range("H45").Value = 1911
$ws.Range("I45").Value = 1833.3334
$ws.Range("K45").Value = 1833.3334
$ws.Range("M45").Value = -1456.3334
$ws.Range("H63").Value = 200002220
$ws.Range("J63").Value = 100001500
$ws.Range("L63").Value = 100001500
$ws.Range("N63").Value = -100002872
$ws.Range("H66").Value = 200002220
$ws.Range("J66").Value = 100001500
$ws.Range("L66").Value = 500007500
$ws.Range("N66").Value = -500014364
$ws.Range("H97").Value = 660.3
$ws.Range("I97").Value = 446
$ws.Range("J97").Value = 922.2222
$ws.Range("K97").Value = 446
$ws.Range("L97").Value = 922.2222
$ws.Range("M97").Value = 50
$ws.Range("N97").Value = -1914.2222
$ws.Range("H122").Value = 10754966
$ws.Range("I122").Value = 13335138
$ws.Range("J122").Value = 4250
$ws.Range("K122").Value = 40005414
$ws.Range("L122").Value = 12750
$ws.Range("M122").Value = -40002964
$ws.Range("N122").Value = -17650

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 4337.6665
$ws.Range("I24").Value = 997.5
$ws.Range("K24").Value = 997.5
$ws.Range("M24").Value = -762.5
$ws.Range("H34").Value = 3000
$ws.Range("J34").Value = 3000
$ws.Range("L34").Value = 3000
$ws.Range("N34").Value = -3228
$ws.Range("H95").Value = 35656
$ws.Range("J95").Value = 35656
$ws.Range("L95").Value = 35656
$ws.Range("N95").Value = -41148
$ws.Range("H107").Value = 31278150
$ws.Range("I107").Value = 17885.072
$ws.Range("K107").Value = 17885.072
$ws.Range("M107").Value = -15965.072

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 530
$ws.Range("I22").Value = 395
$ws.Range("K22").Value = 395
$ws.Range("M22").Value = -45
$ws.Range("H99").Value = 3120
$ws.Range("I99").Value = 2500
$ws.Range("K99").Value = 2500
$ws.Range("M99").Value = -1002
$ws.Range("H122").Value = 2593.2727
$ws.Range("I122").Value = 2603.1
$ws.Range("K122").Value = 7809.299999999999
$ws.Range("M122").Value = -5359.299999999999
$ws.Range("H126").Value = 3120
$ws.Range("I126").Value = 2500
$ws.Range("K126").Value = 7500
$ws.Range("M126").Value = -5030
$ws.Range("H134").Value = 1778.6428
$ws.Range("I134").Value = 1546.2307
$ws.Range("K134").Value = 4638.6921
$ws.Range("M134").Value = -2103.6921

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 621.4666999999999
$ws.Range("I5").Value = 356.72726
$ws.Range("J5").Value = 1349.5
$ws.Range("K5").Value = 1070.18178
$ws.Range("L5").Value = 4048.5
$ws.Range("M5").Value = -958.1817799999999
$ws.Range("N5").Value = -4272.5
$ws.Range("H29").Value = 138
$ws.Range("I29").Value = 160
$ws.Range("K29").Value = 480
$ws.Range("M29").Value = -203
$ws.Range("H129").Value = 3151.7144
$ws.Range("I129").Value = 2307.5
$ws.Range("J129").Value = 4277.3335
$ws.Range("K129").Value = 6922.5
$ws.Range("L129").Value = 12832.0005
$ws.Range("M129").Value = -1922.5
$ws.Range("N129").Value = -22832.0005
$ws.Range("H135").Value = 621.4666999999999
$ws.Range("I135").Value = 356.72726
$ws.Range("J135").Value = 1349.5
$ws.Range("K135").Value = 3210.54534
$ws.Range("L135").Value = 12145.5
$ws.Range("M135").Value = -675.5453400000001
$ws.Range("N135").Value = -17215.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 83334136
$ws.Range("I7").Value = 100000776
$ws.Range("J7").Value = 900
$ws.Range("K7").Value = 100000776
$ws.Range("L7").Value = 900
$ws.Range("M7").Value = -100000664
$ws.Range("N7").Value = -1124
$ws.Range("H40").Value = 2301.923
$ws.Range("I40").Value = 2175
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 2175
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -2039
$ws.Range("N40").Value = -3272
$ws.Range("H46").Value = 1486.3334
$ws.Range("I46").Value = 820.375
$ws.Range("J46").Value = 1699.44
$ws.Range("K46").Value = 820.375
$ws.Range("L46").Value = 1699.44
$ws.Range("M46").Value = -632.375
$ws.Range("N46").Value = -2075.44
$ws.Range("H61").Value = 16794.875
$ws.Range("I61").Value = 16306.556
$ws.Range("K61").Value = 16306.556
$ws.Range("M61").Value = -16104.556
$ws.Range("H113").Value = 16794.875
$ws.Range("I113").Value = 16306.556
$ws.Range("K113").Value = 16306.556
$ws.Range("M113").Value = -14136.556
$ws.Range("H122").Value = 3331.6538
$ws.Range("I122").Value = 2423.6667
$ws.Range("K122").Value = 7271.000100000001
$ws.Range("M122").Value = -4821.000100000001
$ws.Range("H125").Value = 75000
$ws.Range("J125").Value = 75000
$ws.Range("L125").Value = 75000
$ws.Range("N125").Value = -84840
$ws.Range("H126").Value = 83334136
$ws.Range("I126").Value = 100000776
$ws.Range("J126").Value = 900
$ws.Range("K126").Value = 300002328
$ws.Range("L126").Value = 2700
$ws.Range("M126").Value = -299999858
$ws.Range("N126").Value = -7640
$ws.Range("H127").Value = 40000
$ws.Range("J127").Value = 40000
$ws.Range("L127").Value = 40000
$ws.Range("N127").Value = -49920
$ws.Range("H132").Value = 5477
$ws.Range("I132").Value = 2960.611
$ws.Range("J132").Value = 10509.777
$ws.Range("K132").Value = 8881.832999999999
$ws.Range("L132").Value = 31529.331
$ws.Range("M132").Value = -6351.832999999999
$ws.Range("N132").Value = -36589.331

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = $null
$ws.Range("H62").Value = 4767905
$ws.Range("J62").Value = 8666.666999999999
$ws.Range("L62").Value = 8666.666999999999
$ws.Range("N62").Value = -9914.666999999999
$ws.Range("H64").Value = 26999
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 26999
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 26999
$ws.Range("M64").Value = $null
$ws.Range("N64").Value = -27495
$ws.Range("H65").Value = 4767905
$ws.Range("J65").Value = 8666.666999999999
$ws.Range("L65").Value = 43333.335
$ws.Range("N65").Value = -49573.335
$ws.Range("H67").Value = 26999
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 26999
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 26999
$ws.Range("M67").Value = $null
$ws.Range("N67").Value = -28715
$ws.Range("H107").Value = 697.5
$ws.Range("I107").Value = 697.5
$ws.Range("K107").Value = 2092.5
$ws.Range("M107").Value = -172.5
$ws.Range("H122").Value = 1729.037
$ws.Range("I122").Value = 1641.1428
$ws.Range("J122").Value = 2036.6666
$ws.Range("K122").Value = 4923.428400000001
$ws.Range("L122").Value = 6109.9998
$ws.Range("M122").Value = -2473.428400000001
$ws.Range("N122").Value = -11009.9998
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null
